$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Current layout (before edit):
#   row 7  -> CARDIOTON 300MG 20 TABS.
#   row 8  -> CERELAC ...
#   row 9  -> CETAL COLD & FLU 20 CAPLETS
#   row 10 -> FUCICORT CREAM 20 GM
#   row 11 -> PANADOL ADVANCE 500 MG 48 TABLETS
#   row 12 -> VOLTAREN 75MG/3ML 3 AMP.
#   row 13 -> ZURCAL 40MG 14 GASTRO RESISTANT TAB
#   row 14 -> totals
#   row 15 -> footer
#
# New products are inserted alphabetically:
#   LICID LOTION 30 ML        -> between FUCICORT and PANADOL
#   OPLEX-N SYRUP 125ML       -> between FUCICORT and PANADOL (after LICID)
#   XILONE 5MG/5ML SYRUP 100ML-> between VOLTAREN and ZURCAL
#
# Final layout (after edit):
#   row 7  -> CARDIOTON
#   row 8  -> CERELAC
#   row 9  -> CETAL
#   row 10 -> FUCICORT
#   row 11 -> LICID LOTION 30 ML          (new)
#   row 12 -> OPLEX-N SYRUP 125ML         (new)
#   row 13 -> PANADOL
#   row 14 -> VOLTAREN
#   row 15 -> XILONE 5MG/5ML SYRUP 100ML  (new)
#   row 16 -> ZURCAL
#   row 17 -> totals
#   row 18 -> footer
# ---------------------------------------------------------------------------

# 1) Insert two rows for LICID + OPLEX-N right before the current PANADOL row (11).
$ws.Range("A11:Q12").Insert(-4121)
$ws.Range("A7:Q7").Copy($ws.Range("A11:Q11"))
$ws.Range("A7:Q7").Copy($ws.Range("A12:Q12"))

# 2) Insert one row for XILONE right before the current ZURCAL row (now row 16,
#    since PANADOL/VOLTAREN/ZURCAL shifted down by two after the insert above).
$ws.Range("A16:Q16").Insert(-4121)
$ws.Range("A7:Q7").Copy($ws.Range("A16:Q16"))

# ---------------------------------------------------------------------------
# 3) Row heights (match the source report exactly).
# ---------------------------------------------------------------------------
$ws.Rows(7).RowHeight  = 25.5
$ws.Rows(8).RowHeight  = 24.75
$ws.Rows(9).RowHeight  = 25.5
$ws.Rows(10).RowHeight = 24.75
$ws.Rows(11).RowHeight = 25.5
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 24.75
$ws.Rows(14).RowHeight = 25.5
$ws.Rows(15).RowHeight = 24.75
$ws.Rows(16).RowHeight = 25.5
$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 4) Merge the A:B, C:G, H:K, L:M, N:O cell groups for the three brand-new rows
#    (11, 12, 15) — the other data rows already carry their merges forward.
# ---------------------------------------------------------------------------
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

# ---------------------------------------------------------------------------
# 5) Fill in the values for every data row (serial #, name, balance, order
#    limit, price, selling price, transactions).
# ---------------------------------------------------------------------------
function Set-ReportRow($row, $serial, $name, $balance, $orderLimit, $price, $sellPrice, $txns) {
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 8).Value = $balance
    $ws.Cells.Item($row, 12).Value = $orderLimit
    $ws.Cells.Item($row, 14).Value = $price
    $ws.Cells.Item($row, 16).Value = $sellPrice
    $ws.Cells.Item($row, 17).Value = $txns
}

Set-ReportRow 7  1 "CARDIOTON 300MG 20 TABS."             "1:0" "1" "60.00" "30.0000" "0:1"
Set-ReportRow 8  2 "CERELAC رز بدون لبن"                   "0:0" "0" "40.00" "40.0000" "1:0"
Set-ReportRow 9  3 "CETAL COLD & FLU 20 CAPLETS"           "2:1" "1" "36.00" "18.0000" "0:1"
Set-ReportRow 10 4 "FUCICORT CREAM 20 GM"                  "0:0" "1" "70.00" "70.0000" "1:0"
Set-ReportRow 11 5 "LICID LOTION 30 ML"                    "6:0" "1" "40.00" "40.0000" "1:0"
Set-ReportRow 12 6 "OPLEX-N SYRUP 125ML"                   "9:0" "1" "31.00" "31.0000" "1:0"
Set-ReportRow 13 7 "PANADOL ADVANCE 500 MG 48 TABLETS"     "1:3" "1" "92.00" "23.0000" "0:1"
Set-ReportRow 14 8 "VOLTAREN 75MG/3ML 3 AMP."              "5:2" "1" "51.00" "33.6600" "0:2"
Set-ReportRow 15 9 "XILONE 5MG/5ML SYRUP 100ML"            "1:0" "1" "34.00" "34.0000" "1:0"
Set-ReportRow 16 10 "ZURCAL 40MG 14 GASTRO RESISTANT TAB"  "5:0" "1" "96.00" "96.0000" "1:0"

# ---------------------------------------------------------------------------
# 6) Update the totals row (selling-price column sum) and the generated
#    timestamp in the footer row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 16).Value = 415.66000000000002
$ws.Cells.Item(18, 1).Value = "Saturday, 6 September, 2025 10:33 AM"

Write-Output "done"
